# The sheet used to have a small "TARIH / NOT" table in C2:D5. Rebuild it
# as a "Tarih / Not" table in A1:B4 (new header wording/casing, shifted to
# the top-left corner, dates centered in column A, labels centered in B).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old table that lived in C2:D5.
$ws.Range("C2:D5").Clear()

# New header row.
$ws.Range("A1").Value = "Tarih"
$ws.Range("B1").Value = "Not"

# New data rows (dates as date serials, labels as text).
$ws.Range("A2").Value = 46069
$ws.Range("B2").Value = "deneme1"

$ws.Range("A3").Value = 46076
$ws.Range("B3").Value = "DENEME2"

$ws.Range("A4").Value = 46081
$ws.Range("B4").Value = "DENEME3"

# Format the date column like the old date column was formatted.
$ws.Range("A2:A4").NumberFormat = "m/d/yy"

# Center everything in the new table.
$ws.Range("A1:B4").HorizontalAlignment = -4108

# Column widths for the new layout.
$ws.Columns.Item(1).ColumnWidth = 14
$ws.Columns.Item(2).ColumnWidth = 8.26

# Leave the cursor where the author left it.
$ws.Range("H13").Select() | Out-Null
